{"js": "// The document contains a single table of two-digit-by-two-digit\n// multiplication problems (\"NN\u00d7NN=\") laid out five-per-row, with the\n// five \"problem\" rows interleaved with blank rows. This script replaces\n// the problem text in-place, cell by cell, so run/paragraph formatting\n// (font, size, alignment) is preserved exactly.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected at least one table in the document body.\");\n}\n\nconst table = tables.items[0];\n\n// [rowIndex, colIndex, expectedBeforeText, newText]\nconst replacements = [\n  [0, 0, \"84\u00d786=\", \"48\u00d738=\"],\n  [0, 1, \"94\u00d778=\", \"75\u00d736=\"],\n  [0, 2, \"84\u00d786=\", \"11\u00d718=\"],\n  [0, 3, \"74\u00d743=\", \"75\u00d744=\"],\n  [0, 4, \"53\u00d745=\", \"27\u00d798=\"],\n  [4, 0, \"57\u00d784=\", \"73\u00d730=\"],\n  [4, 1, \"91\u00d788=\", \"80\u00d769=\"],\n  [4, 2, \"12\u00d778=\", \"94\u00d742=\"],\n  [4, 3, \"65\u00d788=\", \"82\u00d734=\"],\n  [4, 4, \"91\u00d721=\", \"76\u00d799=\"],\n  [9, 0, \"85\u00d728=\", \"67\u00d764=\"],\n  [9, 1, \"82\u00d796=\", \"72\u00d727=\"],\n  [9, 2, \"54\u00d754=\", \"39\u00d743=\"],\n  [9, 3, \"66\u00d733=\", \"59\u00d745=\"],\n  [9, 4, \"17\u00d769=\", \"23\u00d765=\"],\n  [14, 0, \"57\u00d766=\", \"47\u00d767=\"],\n  [14, 1, \"37\u00d755=\", \"55\u00d784=\"],\n  [14, 2, \"75\u00d791=\", \"40\u00d788=\"],\n  [14, 3, \"91\u00d768=\", \"93\u00d777=\"],\n  [14, 4, \"40\u00d735=\", \"19\u00d717=\"],\n  [19, 0, \"75\u00d784=\", \"91\u00d738=\"],\n  [19, 1, \"72\u00d785=\", \"61\u00d797=\"],\n  [19, 2, \"72\u00d784=\", \"55\u00d793=\"],\n  [19, 3, \"30\u00d743=\", \"28\u00d795=\"],\n  [19, 4, \"65\u00d747=\", \"86\u00d770=\"],\n];\n\n// Search within each target cell individually (rather than a single\n// document-wide search/replace) because several of the original values\n// repeat elsewhere in the table with different replacement targets.\nconst searchResults = [];\nfor (const [row, col, before] of replacements) {\n  const cell = table.getCell(row, col);\n  const found = cell.body.search(before, { matchCase: true, matchWholeWord: true });\n  found.load(\"items\");\n  searchResults.push(found);\n}\nawait context.sync();\n\nsearchResults.forEach((found, i) => {\n  const [, , before, after] = replacements[i];\n  if (found.items.length === 0) {\n    throw new Error(`Could not find \"${before}\" in cell for replacement #${i}.`);\n  }\n  found.items[0].insertText(after, Word.InsertLocation.replace);\n});\nawait context.sync();\n", "ps1": "# The document contains a single table of two-digit-by-two-digit\n# multiplication problems (\"NN\u00d7NN=\") laid out five-per-row, with the\n# five \"problem\" rows interleaved with blank rows. This script replaces\n# the problem text in-place, cell by cell (by 1-based Row/Column), so\n# the existing run/paragraph formatting (font, size, alignment) stays\n# untouched.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$replacements = @(\n  @{row=1;  col=1; before=\"84\u00d786=\"; after=\"48\u00d738=\"},\n  @{row=1;  col=2; before=\"94\u00d778=\"; after=\"75\u00d736=\"},\n  @{row=1;  col=3; before=\"84\u00d786=\"; after=\"11\u00d718=\"},\n  @{row=1;  col=4; before=\"74\u00d743=\"; after=\"75\u00d744=\"},\n  @{row=1;  col=5; before=\"53\u00d745=\"; after=\"27\u00d798=\"},\n  @{row=5;  col=1; before=\"57\u00d784=\"; after=\"73\u00d730=\"},\n  @{row=5;  col=2; before=\"91\u00d788=\"; after=\"80\u00d769=\"},\n  @{row=5;  col=3; before=\"12\u00d778=\"; after=\"94\u00d742=\"},\n  @{row=5;  col=4; before=\"65\u00d788=\"; after=\"82\u00d734=\"},\n  @{row=5;  col=5; before=\"91\u00d721=\"; after=\"76\u00d799=\"},\n  @{row=10; col=1; before=\"85\u00d728=\"; after=\"67\u00d764=\"},\n  @{row=10; col=2; before=\"82\u00d796=\"; after=\"72\u00d727=\"},\n  @{row=10; col=3; before=\"54\u00d754=\"; after=\"39\u00d743=\"},\n  @{row=10; col=4; before=\"66\u00d733=\"; after=\"59\u00d745=\"},\n  @{row=10; col=5; before=\"17\u00d769=\"; after=\"23\u00d765=\"},\n  @{row=15; col=1; before=\"57\u00d766=\"; after=\"47\u00d767=\"},\n  @{row=15; col=2; before=\"37\u00d755=\"; after=\"55\u00d784=\"},\n  @{row=15; col=3; before=\"75\u00d791=\"; after=\"40\u00d788=\"},\n  @{row=15; col=4; before=\"91\u00d768=\"; after=\"93\u00d777=\"},\n  @{row=15; col=5; before=\"40\u00d735=\"; after=\"19\u00d717=\"},\n  @{row=20; col=1; before=\"75\u00d784=\"; after=\"91\u00d738=\"},\n  @{row=20; col=2; before=\"72\u00d785=\"; after=\"61\u00d797=\"},\n  @{row=20; col=3; before=\"72\u00d784=\"; after=\"55\u00d793=\"},\n  @{row=20; col=4; before=\"30\u00d743=\"; after=\"28\u00d795=\"},\n  @{row=20; col=5; before=\"65\u00d747=\"; after=\"86\u00d770=\"}\n)\n\nforeach ($rep in $replacements) {\n  $cell = $t.Cell($rep.row, $rep.col)\n  $range = $cell.Range\n  # Cell ranges carry a trailing cell-mark (CR + BEL); strip it before comparing.\n  $current = $range.Text.TrimEnd([char]13, [char]7)\n  if ($current -ne $rep.before) {\n    throw \"Cell ($($rep.row),$($rep.col)) expected '$($rep.before)' but found '$current'\"\n  }\n  $range.Text = $rep.after\n}\n"}
